$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header summary cells
# ---------------------------------------------------------------------------
$ws.Range("E11").Value = 644707
$ws.Range("C13").Value = 4
$ws.Range("F13").Value = 11

# ---------------------------------------------------------------------------
# Before touching the detail rows, copy the special "last row" border
# formatting (currently on row 32) onto row 26, which will become the new
# last data row once the surplus rows are removed below.
# ---------------------------------------------------------------------------
$ws.Range("B32:J32").Copy()
$ws.Range("B26:J26").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Rewrite the worker detail table (rows 16-26)
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "1001941167"
$ws.Range("D16").Value = "ERIKA PAOLA DE MOYA PARDO"
$ws.Range("E16").Value = "1903"
$ws.Range("F16").Value = 63400
$ws.Range("G16").Value = 1902000

$ws.Range("C17").Value = "1140909509"
$ws.Range("D17").Value = "WILMER JOSE PAZ CARRASQUERO"
$ws.Range("E17").Value = "2001"
$ws.Range("F17").Value = 39258
$ws.Range("G17").Value = 2859276

$ws.Range("C18").Value = "1140909509"
$ws.Range("D18").Value = "WILMER JOSE PAZ CARRASQUERO"
$ws.Range("E18").Value = "2002"
$ws.Range("F18").Value = 73609
$ws.Range("G18").Value = 2859276

$ws.Range("C19").Value = "1140909509"
$ws.Range("D19").Value = "WILMER JOSE PAZ CARRASQUERO"
$ws.Range("E19").Value = "2003"
$ws.Range("F19").Value = 73609
$ws.Range("G19").Value = 2859276

$ws.Range("C20").Value = "1140909509"
$ws.Range("D20").Value = "WILMER JOSE PAZ CARRASQUERO"
$ws.Range("E20").Value = "2004"
$ws.Range("F20").Value = 73609
$ws.Range("G20").Value = 2859276

$ws.Range("C21").Value = "1140909509"
$ws.Range("D21").Value = "WILMER JOSE PAZ CARRASQUERO"
$ws.Range("E21").Value = "2005"
$ws.Range("F21").Value = 73609
$ws.Range("G21").Value = 2859276

$ws.Range("C22").Value = "1140909509"
$ws.Range("D22").Value = "WILMER JOSE PAZ CARRASQUERO"
$ws.Range("E22").Value = "2006"
$ws.Range("F22").Value = 73609
$ws.Range("G22").Value = 2859276

$ws.Range("C23").Value = "1140909509"
$ws.Range("D23").Value = "WILMER JOSE PAZ CARRASQUERO"
$ws.Range("E23").Value = "2007"
$ws.Range("F23").Value = 73609
$ws.Range("G23").Value = 2859276

$ws.Range("C24").Value = "1140909509"
$ws.Range("D24").Value = "WILMER JOSE PAZ CARRASQUERO"
$ws.Range("E24").Value = "2008"
$ws.Range("F24").Value = 73609
$ws.Range("G24").Value = 2859276

$ws.Range("C25").Value = "1001970826"
$ws.Range("D25").Value = "CAROLINA ISABEL DIAZ BOSSIO"
$ws.Range("E25").Value = "2102"
$ws.Range("F25").Value = 16386
$ws.Range("G25").Value = 1228924

$ws.Range("C26").Value = "1143346433"
$ws.Range("D26").Value = "ALEXANDER MANUEL MORENO BABILONIA"
$ws.Range("E26").Value = "2401"
$ws.Range("F26").Value = 10400
$ws.Range("G26").Value = 1300000

# ---------------------------------------------------------------------------
# Remove the now-obsolete rows (old rows 27-32), shifting the footer rows
# (37/38 -> 31/32) up automatically.
# ---------------------------------------------------------------------------
$ws.Range("B27:J32").Delete(-4162)

# ---------------------------------------------------------------------------
# Column D is now narrower because the longest worker name was removed;
# update the stored best-fit width to match the new longest entry.
# ---------------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 38.36328125
